$d = $word.ActiveDocument

$replacements = @(
    @("466×4=", "940×6="),
    @("801×5=", "915×7="),
    @("683×4=", "979×9="),
    @("740×3=", "779×5="),
    @("523×4=", "904×8="),
    @("674×5=", "705×6="),
    @("389×7=", "346×7="),
    @("696×6=", "697×7="),
    @("197×9=", "125×2="),
    @("269×9=", "268×4="),
    @("389×2=", "728×9="),
    @("484×3=", "486×2="),
    @("105×2=", "977×6="),
    @("517×3=", "913×4="),
    @("327×7=", "147×2="),
    @("481×7=", "181×7="),
    @("422×4=", "151×2="),
    @("331×8=", "260×4="),
    @("741×5=", "990×4="),
    @("701×4=", "893×7="),
    @("220×5=", "410×2="),
    @("981×7=", "317×6="),
    @("244×6=", "257×8="),
    @("299×2=", "521×7="),
    @("817×9=", "220×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
